$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.810.52'
$ws.Range("E2").Value = '  -1.13%  '

$ws.Range("D3").Value = '3.847.18'
$ws.Range("E3").Value = '  -1.59%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.59'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.94'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("D7").Value = '3.845.18'
$ws.Range("E7").Value = '  -1.65%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("E10").Value = '  -0.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.34'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -0.17%  '

$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.75'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("D15").Value = '4.495.78'
$ws.Range("E15").Value = '  -1.37%  '

$ws.Range("D16").Value = '3.856.66'
$ws.Range("E16").Value = '  -1.38%  '

$ws.Range("D17").Value = '67.919.95'
$ws.Range("E17").Value = '  -1.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.06'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +6.06%  '

$ws.Range("E19").Value = '  -1.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.111'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -1.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.92'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '462.47'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -4.50%  '

$ws.Range("E23").Value = '  +1.26%  '

$ws.Range("E24").Value = '  -3.81%  '

$ws.Range("E25").Value = '  -1.52%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.08'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.94'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -1.32%  '

$ws.Range("E30").Value = '  +0.30%  '

$ws.Range("D31").Value = '3.998.59'
$ws.Range("E31").Value = '  -1.42%  '

$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("E33").Value = '  -2.39%  '

$ws.Range("E34").Value = '  -3.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.27'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +1.79%  '

$ws.Range("D36").Value = '3.824.86'
$ws.Range("E36").Value = '  -0.63%  '

$ws.Range("E37").Value = '  -1.95%  '

$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("E40").Value = '  -0.20%  '

$ws.Range("E41").Value = '  +5.61%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.310'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -1.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '424.78'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -1.57%  '

$ws.Range("E45").Value = '  -0.47%  '

$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("E48").Value = '  +0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000274'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +4.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '143.69'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.25%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.41'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.42%  '
